# ADD completed the two tables and finished the CRUD V1.3
# Rewrites the sample data rows (12-19) of the details table and removes the
# now-redundant last row (20), collapsing the used range from A1:H20 to A1:H19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 --------------------------------------------------------------
$ws.Range("A12").Value = "4857856"
$ws.Range("B12").Value = "PLANTILLA DE APORTES"
$ws.Range("C12").Value = "32"
$ws.Range("D12").Value = "23"
$ws.Range("E12").Value = "PLANTILLA DE APORTES"
$ws.Range("F12").Value = 32
$ws.Range("G12").Value = 23
$ws.Range("H12").Value = "Dr. Atiencia Atiencia Atiencia Atiencia"

# --- Row 13 --------------------------------------------------------------
$ws.Range("A13").Value = "45785869"
$ws.Range("B13").Value = "Raul Alejandro Sosa"
$ws.Range("C13").Value = "172845688978"
$ws.Range("D13").Value = "1548785225"
$ws.Range("E13").Value = "PLANTILLA DE APORTES"
$ws.Range("F13").Value = 45.32
$ws.Range("G13").Value = 25.26
$ws.Range("H13").Value = "Dr. Christian Santiago Izurieta Cruz"

# --- Row 14 --------------------------------------------------------------
$ws.Range("A14").Value = "45785887"
$ws.Range("B14").Value = "Raul Alejandro Sosa"
$ws.Range("C14").Value = "172845688978"
$ws.Range("D14").Value = "1548785225"
$ws.Range("E14").Value = "PLANTILLA DE APORTES"
$ws.Range("F14").Value = 45.25
$ws.Range("G14").Value = 15.26
$ws.Range("H14").Value = "Dr. Christian Santiago Izurieta Cruz"

# --- Row 15 --------------------------------------------------------------
$ws.Range("A15").Value = "458965"
$ws.Range("B15").Value = "Alexander Benitez"
$ws.Range("C15").Value = "1458789"
$ws.Range("D15").Value = "1728224557"
$ws.Range("E15").Value = "FONDO DE RESERVA"
$ws.Range("F15").Value = 45.23
$ws.Range("G15").Value = 12.25
$ws.Range("H15").Value = "Dr. Atiencia Atiencia Atiencia Atiencia"

# --- Row 16 --------------------------------------------------------------
$ws.Range("A16").Value = "235689"
$ws.Range("B16").Value = "Alexander Benitez"
$ws.Range("C16").Value = "1458789"
$ws.Range("D16").Value = "1728224557"
$ws.Range("E16").Value = "PLANTILLA DE APORTES"
$ws.Range("F16").Value = 45.23
$ws.Range("G16").Value = 12.25
$ws.Range("H16").Value = "Dr. Atiencia Atiencia Atiencia Atiencia"

# --- Row 17 --------------------------------------------------------------
$ws.Range("A17").Value = "457896"
$ws.Range("B17").Value = "Josue Alberto Ramirez Arboleda"
$ws.Range("C17").Value = "174578569933"
$ws.Range("D17").Value = "1245785689"
$ws.Range("E17").Value = "PLANTILLA DE APORTES"
$ws.Range("F17").Value = 12.25
$ws.Range("G17").Value = 12.23
$ws.Range("H17").Value = "Dr. Christian Santiago Izurieta Cruz"

# --- Row 18 --------------------------------------------------------------
$ws.Range("A18").Value = "451278"
$ws.Range("B18").Value = "Josue Alberto Ramirez Arboleda"
$ws.Range("C18").Value = "174578569933"
$ws.Range("D18").Value = "1245785689"
$ws.Range("E18").Value = "PLANTILLA DE APORTES"
$ws.Range("F18").Value = 150.26
$ws.Range("G18").Value = 325.26
$ws.Range("H18").Value = "Dr. Christian Santiago Izurieta Cruz"

# --- Row 19 --------------------------------------------------------------
$ws.Range("A19").Value = "561245"
$ws.Range("B19").Value = "Josue Alberto Ramirez Arboleda"
$ws.Range("C19").Value = "174578569933"
$ws.Range("D19").Value = "1245785689"
$ws.Range("E19").Value = "PRESTAMOS"
$ws.Range("F19").Value = 152.23
$ws.Range("G19").Value = 321.23
$ws.Range("H19").Value = "Dr. Christian Santiago Izurieta Cruz"

# --- Row 20 is now redundant (its data was folded into row 12) ----------
$ws.Rows(20).Delete()
